# Weekly data refresh: a new price observation is published for
# "Terminal La Palmera de La Serena - Zanahoria", so a new row is
# inserted at row 503 (pushing the existing 503:571 rows down to
# 504:572) and populated with the newest record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 503; this shifts rows 503:571 down to 504:572
# and keeps the sheet's used-range/dimension in sync automatically.
$ws.Rows("503:503").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A503").Value = 8
$ws.Range("B503").Value = "Terminal La Palmera de La Serena"
$ws.Range("C503").Value = "Coquimbo"
$ws.Range("D503").Value = 45131
$ws.Range("E503").Value = 4
$ws.Range("F503").Value = 100114013
$ws.Range("G503").Value = "Zanahoria"
$ws.Range("H503").Value = "Sin especificar"
$ws.Range("I503").Value = "Primera"
$ws.Range("J503").Value = 600
$ws.Range("K503").Value = 5800
$ws.Range("L503").Value = 6000
$ws.Range("M503").Value = 5900
$ws.Range("N503").Value = "$/saco 20 kilos"
$ws.Range("O503").Value = "Provincia del Elquí"
$ws.Range("P503").Value = 295
$ws.Range("Q503").Value = 20
$ws.Range("R503").Value = "Hortaliza"
